# Apply the style changes described by the commit:
#  1. Add a new paragraph style "AbstractTitle" ("Abstract Title"),
#     based on Normal, followed by Abstract, centered/bold/colored.
#  2. Change the "Abstract" style's space-before from 300 (15pt) to
#     100 (5pt) twentieths-of-a-point (space-after stays 300/15pt).

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" paragraph style -----------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)   # 1 = wdStyleTypeParagraph
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true     # <w:keepNext/>
$abstractTitle.ParagraphFormat.KeepTogether = $true      # <w:keepLines/>
$abstractTitle.ParagraphFormat.Alignment = 1              # wdAlignParagraphCenter -> <w:jc w:val="center"/>
$abstractTitle.ParagraphFormat.SpaceBefore = 15            # 300 twips -> <w:spacing w:before="300" .../>
$abstractTitle.ParagraphFormat.SpaceAfter = 0               # <w:spacing .../w:after="0"/>

$abstractTitle.Font.Size = 10      # <w:sz w:val="20"/>
$abstractTitle.Font.SizeBi = 10     # <w:szCs w:val="20"/>
$abstractTitle.Font.Bold = $true     # <w:b/>
$abstractTitle.Font.Color = 9067060   # 0x8A5A34 (BGR) == RGB 345A8A -> <w:color w:val="345A8A"/>

# --- 2. "Abstract" style spacing tweak ------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5   # 100 twips -> <w:spacing w:before="100" .../>

Write-Output "AbstractTitle style added; Abstract spacing updated."
